$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Number column values for rows 5, 6, 7
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5

# Update the active selection to B3:B7 with active cell B3
$ws.Range("B3:B7").Select()
